{"js": "// Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line,\n// the \"\u00a9 2020 ... Creative Commons Attribution\" footer line, and the blank\n// paragraph that separates them from the \"LOM3213: Fen\u00f4menos de Transporte B\n// (Requisito)\" requirement line above them.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its text so the edit is resilient to the\n// exact paragraph index.\nconst anchor = paragraphs.items.find(\n  (p) => p.text.trim() === \"LOM3213: Fen\u00f4menos de Transporte B (Requisito)\"\n);\nif (!anchor) {\n  throw new Error('Anchor paragraph \"LOM3213: ...\" not found.');\n}\n\n// The three paragraphs immediately following the anchor are:\n//   1. an empty paragraph\n//   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3. \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//       pages. Original theme under Creative Commons Attribution\"\n// All three should be removed, leaving the anchor paragraph followed\n// directly by the paragraph that used to come after them.\nconst toDelete = [];\nlet current = anchor;\nfor (let i = 0; i < 3; i++) {\n  current = current.getNext();\n  current.load(\"text\");\n  await context.sync();\n  toDelete.push(current);\n}\n\nconst expectedTexts = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\ntoDelete.forEach((p, i) => {\n  if (p.text !== expectedTexts[i]) {\n    throw new Error(\n      `Unexpected paragraph text while deleting footer block: got \"${p.text}\", expected \"${expectedTexts[i]}\"`\n    );\n  }\n});\n\n// Delete from the last one backwards so each deleted paragraph's position\n// stays valid during the loop.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line,\n# the \"\u00a9 2020 ... Creative Commons Attribution\" footer line, and the blank\n# paragraph that separates them from the \"LOM3213: Fen\u00f4menos de Transporte B\n# (Requisito)\" requirement line above them.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its text so the edit is resilient to the\n# exact paragraph index.\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd() -eq \"LOM3213: Fen\u00f4menos de Transporte B (Requisito)\") {\n        $anchor = $p\n        break\n    }\n}\nif ($anchor -eq $null) {\n    throw \"Anchor paragraph 'LOM3213: ...' not found.\"\n}\n\n# The three paragraphs immediately following the anchor are:\n#   1. an empty paragraph\n#   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3. \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#       pages. Original theme under Creative Commons Attribution\"\n# All three should be removed, leaving the anchor paragraph followed\n# directly by the paragraph that used to come after them.\n$p1 = $anchor.Next()\n$p2 = $p1.Next()\n$p3 = $p2.Next()\n\n$expected1 = \"\"\n$expected2 = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$expected3 = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\nif ($p1.Range.Text.TrimEnd() -ne $expected1) {\n    throw \"Unexpected paragraph text while deleting footer block (1).\"\n}\nif ($p2.Range.Text.TrimEnd() -ne $expected2) {\n    throw \"Unexpected paragraph text while deleting footer block (2).\"\n}\nif ($p3.Range.Text.TrimEnd() -ne $expected3) {\n    throw \"Unexpected paragraph text while deleting footer block (3).\"\n}\n\n# Delete from the last one backwards so each deleted paragraph's position\n# stays valid during the operation.\n$p3.Range.Delete()\n$p2.Range.Delete()\n$p1.Range.Delete()\n"}
